# Auto-generated: applies cryptos list price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.554.69"
$ws.Range("E2").Value = "  +2.90%  "
$ws.Range("D3").Value = "1.786.69"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'224.30"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "'0.557"
$ws.Range("E6").Value = "  +1.37%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "'33.02"
$ws.Range("E8").Value = "  +9.32%  "
$ws.Range("D9").Value = "'0.282"
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("D10").Value = "'0.0680"
$ws.Range("E10").Value = "  +3.84%  "
$ws.Range("D11").Value = "'0.0935"
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("D12").Value = "2.046.66"
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").Value = "'11.04"
$ws.Range("E13").Value = "  +11.77%  "
$ws.Range("D14").Value = "1.793.04"
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").Value = "'0.633"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").Value = "34.552.97"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("D17").Value = "'4.28"
$ws.Range("E17").Value = "  +3.18%  "
$ws.Range("D18").Value = "'68.59"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").Value = "'253.69"
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("D20").Value = "0.0₃0774"
$ws.Range("E20").Value = "  +5.83%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "'10.42"
$ws.Range("E22").Value = "  +2.45%  "
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").Value = "'2.14"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").Value = "'158.62"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").Value = "'16.35"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "'7.09"
$ws.Range("E27").Value = "  +3.54%  "
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  +1.56%  "
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").Value = "'3.58"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("E34").Value = "  +4.90%  "
$ws.Range("D35").Value = "1.443.03"
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("E37").Value = "  +3.58%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").Value = "'83.15"
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("D40").Value = "'2.80"
$ws.Range("E40").Value = "  +4.52%  "
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("D42").Value = "'0.895"
$ws.Range("E42").Value = "  +2.15%  "
$ws.Range("E43").Value = "  +0.76%  "
$ws.Range("D44").Value = "'0.0503"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D45").Value = "'5.89"
$ws.Range("E45").Value = "  +2.84%  "
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("D47").Value = "1.942.81"
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("D48").Value = "'104.45"
$ws.Range("E48").Value = "  +7.76%  "
$ws.Range("D49").Value = "'12.00"
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").Value = "0.0₆0121"
$ws.Range("E51").Value = "  +6.58%  "
